$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell F1: "Material Type" ---
# Start from the formatting of an existing bold/centered header cell (E1) so that
# the bold font (fontId=1, theme white) and alignment/fill "shape" are reused
# instead of building a brand new font+alignment combo from scratch.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false

$h = $ws.Range("F1")
$h.Value = "Material Type"
$h.Font.Color = 16777215       # explicit white (was theme-based white)
$h.Interior.Color = 0          # black fill foreground
$h.Interior.PatternColor = 0   # black fill background

# --- Data cells F2:F24: "DNA:Genomic", centered (reuses existing style 3) ---
for ($row = 2; $row -le 24; $row++) {
    $cell = $ws.Cells.Item($row, 6)
    $cell.Value = "DNA:Genomic"
    $cell.HorizontalAlignment = -4108  # xlCenter
}

# Update the active selection to reflect the newly added column, as in the target sheet
$ws.Range("F1:F24").Select()
